$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns hold text-formatted values (e.g. "1.003",
# "0.0₅8129", "  -1.31%  "). Force Text number format first so Excel's COM
# Range.Value setter doesn't coerce numeric-looking strings into real numbers
# (which would e.g. turn "0.4250" into 0.425 or "1.003" into 1.0029999999999999).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.204.78'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.659.99'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.21'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5183'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.00%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06279'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07765'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.485'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.661.94'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.886.26'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5474'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8129'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.02'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.209.69'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.617'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.38'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.08'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.016'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -4.84%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.48'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.297'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.16'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.08%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05938'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.275'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.549'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.284'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.584'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -6.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9621'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.418'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.769'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5674'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -6.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.037'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8529'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.012.56'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -7.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.96'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.801.04'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈110'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.54'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.30%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.025'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.25%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4250'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.55%  '
